# Updated cryptos list on Fri Mar  3 07:48:26 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.409.87'
$ws.Range("E2").Value = '  -4.44%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.571.34'
$ws.Range("E3").Value = '  -4.59%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.000'
$ws.Range("E5").Value = '  -0.14%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '291.29'
$ws.Range("E6").Value = '  -2.67%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3676'
$ws.Range("E7").Value = '  -3.00%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.31'
$ws.Range("E8").Value = '  -0.92%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3375'
$ws.Range("E9").Value = '  -5.17%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.173'
$ws.Range("E10").Value = '  -3.83%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07584'
$ws.Range("E11").Value = '  -6.34%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9996'
$ws.Range("E12").Value = '  -0.19%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.10'
$ws.Range("E13").Value = '  -4.19%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.057'
$ws.Range("E14").Value = '  -5.18%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.874'
$ws.Range("E15").Value = '  -6.65%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001144'
$ws.Range("E16").Value = '  -4.26%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.568.32'
$ws.Range("E17").Value = '  -4.73%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '89.14'
$ws.Range("E18").Value = '  -8.38%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06721'
$ws.Range("E19").Value = '  -3.31%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  -0.12%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.260'
$ws.Range("E21").Value = '  -7.47%  '

$ws.Range("E22").Value = '  -5.20%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.5243'
$ws.Range("E23").Value = '  -8.91%  '

$ws.Range("E24").Value = '  -3.31%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '22.414.30'
$ws.Range("E25").Value = '  -4.50%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.383'
$ws.Range("E26").Value = '  -4.57%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.984'
$ws.Range("E27").Value = '  +1.85%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.87'
$ws.Range("E28").Value = '  -4.79%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '145.43'
$ws.Range("E29").Value = '  -4.84%  '

$ws.Range("E30").Value = '  -5.02%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '124.98'
$ws.Range("E31").Value = '  -5.78%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.746.90'
$ws.Range("E32").Value = '  -4.51%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.279'
$ws.Range("E33").Value = '  -9.23%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.007'
$ws.Range("E34").Value = '  -0.22%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.974'
$ws.Range("E35").Value = '  -6.02%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.38'
$ws.Range("E36").Value = '  -11.64%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.08432'

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02542'
$ws.Range("E38").Value = '  -6.61%  '

$ws.Range("E39").Value = '  -5.33%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.535'
$ws.Range("E40").Value = '  -6.55%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.06517'
$ws.Range("E41").Value = '  -3.81%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.80'
$ws.Range("E42").Value = '  -9.60%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.252'
$ws.Range("E43").Value = '  -3.77%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6392'
$ws.Range("E44").Value = '  -7.09%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.60'
$ws.Range("E45").Value = '  -6.28%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9995'
$ws.Range("E46").Value = '  -0.19%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6030'
$ws.Range("E47").Value = '  -5.69%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.764'
$ws.Range("E48").Value = '  -3.86%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.122'

$ws.Range("B50").Value = 'EOS'
$ws.Range("C50").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.200'
$ws.Range("E50").Value = '  +2.29%  '

$ws.Range("B51").Value = 'Quant'
$ws.Range("C51").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '121.22'
$ws.Range("E51").Value = '  -4.83%  '
